$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Fix the Expected_result (column K) csv paths for rows 117-122: the path
#    segment "testdata/cases/..." should read "testdata/mysqlcases/...".
# ---------------------------------------------------------------------------
$ws.Range("K117").Value = "src/test/resources/io.dingodb.test/testdata/mysqlcases/dml/updatedelete/expectedresult/updatedelete_116.csv"
$ws.Range("K118").Value = "src/test/resources/io.dingodb.test/testdata/mysqlcases/dml/updatedelete/expectedresult/updatedelete_117.csv"
$ws.Range("K119").Value = "src/test/resources/io.dingodb.test/testdata/mysqlcases/dml/updatedelete/expectedresult/updatedelete_118.csv"
$ws.Range("K120").Value = "src/test/resources/io.dingodb.test/testdata/mysqlcases/dml/updatedelete/expectedresult/updatedelete_119.csv"
$ws.Range("K121").Value = "src/test/resources/io.dingodb.test/testdata/mysqlcases/dml/updatedelete/expectedresult/updatedelete_120.csv"
$ws.Range("K122").Value = "src/test/resources/io.dingodb.test/testdata/mysqlcases/dml/updatedelete/expectedresult/updatedelete_121.csv"

# ---------------------------------------------------------------------------
# 2) Append two brand new test rows (123 & 124) for Chinese-character data
#    cases. Inserting rows (rather than just writing past the used range)
#    makes the new rows inherit row 122's formatting/style.
# ---------------------------------------------------------------------------
$ws.Rows("123:124").Insert()

# Row 123: updel_122 - update Chinese field value
$ws.Range("A123").Value = "updel_122"
$ws.Range("B123").Value = "y"
$ws.Range("C123").Value = "更新中文字段值"
$ws.Range("D123").Value = "SQLFunction"
$ws.Range("F123").Value = "schema1"
$ws.Range("G123").Value = "updel_value05"
$ws.Range("H123").Value = "update `$schema1 set name='眼前无路想回头' where id=1 or id=2"
$ws.Range("I123").Value = "2"
$ws.Range("J123").Value = "select id,name,age,amount,address,birthday,create_time,update_time,is_delete from `$schema1"
$ws.Range("K123").Value = "src/test/resources/io.dingodb.test/testdata/mysqlcases/dml/updatedelete/expectedresult/updatedelete_122.csv"
$ws.Range("L123").Value = "csv_containsAll"

# Row 124: updel_123 - delete Chinese data
$ws.Range("A124").Value = "updel_123"
$ws.Range("B124").Value = "y"
$ws.Range("C124").Value = "删除中文数据"
$ws.Range("D124").Value = "SQLFunction"
$ws.Range("F124").Value = "schema1"
$ws.Range("G124").Value = "updel_value05"
$ws.Range("H124").Value = "delete from `$schema1 where address like '%测试%' or address='北京纯牛奶'"
$ws.Range("I124").Value = "2"
$ws.Range("J124").Value = "select * from `$schema1 where address='北京纯牛奶' or address like '%测试%' or address in ('上海虹桥')"
$ws.Range("K124").Value = "src/test/resources/io.dingodb.test/testdata/mysqlcases/dml/updatedelete/expectedresult/updatedelete_123.csv"
$ws.Range("L124").Value = "csv_containsAll"

# ---------------------------------------------------------------------------
# 3) A handful of cells in the two new rows use the plain (unshaded) cell
#    style rather than the shaded one inherited from row 122 - clear their
#    fill to match.
# ---------------------------------------------------------------------------
foreach ($addr in @("C123","G123","H123","I123","J123","C124","H124","I124")) {
    $ws.Range($addr).Interior.Pattern = -4142
}

# ---------------------------------------------------------------------------
# 4) Restore the active-cell selection shown in the sheet view.
# ---------------------------------------------------------------------------
$ws.Range("B95").Select()
